$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell M1: "valorRecarga" -> "valorDescarga"
$ws.Range("M1").Value = "valorDescarga"

# Move the active selection to M1 (matches the edited cell)
$ws.Range("M1").Select()
